# Generate Report for Handoff
# Replace the old GUID-based file identifiers and refresh the handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "2bdc9880-ef42-4f49-8015-389b96571506"
$newGuid = "3e1e5e9f-8674-4c80-93ce-b59daefaedb5"

$oldHash = "62da0c859296ee3475ff38a5ebb49097dde1948f"
$newHash = "f34def58fbd07448192c3e5344ef55fd5b51a18a"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-16 04:54:57"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"

# --- zh-cn sheet ---
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-16 04:54:52"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"

# --- de-de sheet ---
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-16 04:54:57"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
